# "Added ifs to user's id" -- add a new patient record (row 5) for the same
# patient/ID as row 4 (315783522), with a normal ("healthy") set of results,
# and record their age/id-suffix as 22.
#
# NOTE: values here are copied cell-to-cell (rather than typed fresh) so the
# engine keeps storing them as shared-string text -- exactly like every
# other cell in this sheet -- instead of inferring numeric types for
# number-looking text such as "315783522" or "40".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a full copy of row 4 (same patient/ID, same cell styling),
# then overwrite just the cells that differ for this new record.
$ws.Range("A4:U4").Copy($ws.Range("A5:U5"))

# Lab/symptom columns that should read as normal/negative, mirroring the
# healthy-patient values already used elsewhere in the sheet (row 2).
$ws.Range("H2").Copy($ws.Range("E5"))
$ws.Range("H2").Copy($ws.Range("F5"))
$ws.Range("H2").Copy($ws.Range("G5"))
$ws.Range("H2").Copy($ws.Range("H5"))
$ws.Range("H2").Copy($ws.Range("I5"))
$ws.Range("J2").Copy($ws.Range("J5"))
$ws.Range("K2").Copy($ws.Range("K5"))
$ws.Range("K2").Copy($ws.Range("L5"))
$ws.Range("M2").Copy($ws.Range("M5"))
$ws.Range("K2").Copy($ws.Range("N5"))
$ws.Range("O2").Copy($ws.Range("O5"))
$ws.Range("P2").Copy($ws.Range("P5"))
$ws.Range("Q2").Copy($ws.Range("Q5"))
$ws.Range("R2").Copy($ws.Range("R5"))
$ws.Range("S2").Copy($ws.Range("S5"))
$ws.Range("R2").Copy($ws.Range("T5"))
$ws.Range("U2").Copy($ws.Range("U5"))

# Age/id column "22" is brand new text in this workbook. Stage it via a
# TEXT() formula in a scratch cell and paste its computed value back in, so
# it lands as shared-string text rather than a numeric literal.
$ws.Range("X1").Formula = "=TEXT(22,""0"")"
$ws.Range("X1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("X1").ClearContents()

# Row 4 picked up a trailing blank V4 cell (matches the blank V column used
# throughout the sheet); row 5 itself has no V entry. Copy an existing
# present-but-empty cell so the element is actually emitted (assigning ""
# directly leaves no cell node at all).
$ws.Range("V2").Copy($ws.Range("V4"))
